$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = -1.749529265654429
$ws.Range("B15").Value = 0.8615440086535034
$ws.Range("B16").Value = -0.8699571331276832
$ws.Range("B17").Value = -0.8699571331276833
$ws.Range("B18").Value = -0.8795721325267457
$ws.Range("B19").Value = -1.000000000000005
$ws.Range("B23").Value = 3.604823524698546
$ws.Range("B25").Value = 1.855294259044117
$ws.Range("B26").Value = 1.855294259044117
$ws.Range("B27").Value = 1.855294259044118
$ws.Range("B28").Value = 2.793037137935192
$ws.Range("B29").Value = -2.793037137935192
$ws.Range("B30").Value = -2.761548014903262
$ws.Range("B31").Value = 2.761548014903262
$ws.Range("B32").Value = 1.705781018388694
$ws.Range("B33").Value = -0.03365249789671901
$ws.Range("B40").Value = 0.006970874564320367
$ws.Range("B42").Value = 0.4106806618324632
$ws.Range("B44").Value = 0.4003445374784709
$ws.Range("B45").Value = 1.534313529105413
$ws.Range("B46").Value = 0
$ws.Range("B72").Value = 0
$ws.Range("B77").Value = -0.008893874444132883
$ws.Range("B78").Value = 0
$ws.Range("B79").Value = 0
$ws.Range("B80").Value = 0
$ws.Range("B81").Value = 0
$ws.Range("B119").Value = -2.844116822242711
$ws.Range("B120").Value = -0.7966027002123338
$ws.Range("B121").Value = 0.7966027002123338
$ws.Range("B122").Value = 4.192500300468752
$ws.Range("B127").Value = 1.060414246224113
$ws.Range("B129").Value = 1.060414246224113
$ws.Range("B130").Value = 0.008893874444132883
$ws.Range("B136").Value = [double]"2.881211122307894e-16"
$ws.Range("B137").Value = 5.688233644485422
$ws.Range("B138").Value = -0.2072032370497985
$ws.Range("B139").Value = 0.2072032370497985
$ws.Range("B140").Value = -0.005047874684507852
$ws.Range("B141").Value = 0.005047874684507852
$ws.Range("B143").Value = 0.1139377428788915
$ws.Range("B144").Value = -0.02067224870798434
$ws.Range("B145").Value = 0.02067224870798434
$ws.Range("B147").Value = 0.02403749849765644
$ws.Range("B148").Value = 0.02403749849765644
$ws.Range("B151").Value = 0.006970874564320367
$ws.Range("B153").Value = 0.4106806618324632
$ws.Range("B155").Value = 0
$ws.Range("B170").Value = 1.032290372981856
$ws.Range("B173").Value = 1.032290372981856
$ws.Range("B175").Value = 2.742558391090113
$ws.Range("B188").Value = 3.709466768158342
$ws.Range("B192").Value = -0.8653499459156319
$ws.Range("B196").Value = -1.032290372981856
$ws.Range("B198").Value = 1.799447137534562
$ws.Range("B200").Value = -1.101037618685153
$ws.Range("B201").Value = -0.7671567645527064
$ws.Range("B209").Value = 0
$ws.Range("B210").Value = 1.000000000000005
$ws.Range("B215").Value = 0
$ws.Range("B226").Value = -2.611073274307932
$ws.Range("B227").Value = -2.611073274307932
$ws.Range("B228").Value = -2.611073274307932
$ws.Range("B229").Value = -2.611073274307932
